$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-19 Wednesday" "2025-02-20 Thursday"

Replace-Text "79÷7=" "32÷4="
Replace-Text "82÷4=" "95÷7="
Replace-Text "56÷9=" "58÷9="
Replace-Text "84÷7=" "95÷3="
Replace-Text "94÷5=" "45÷6="

Replace-Text "57÷7=" "78÷9="
Replace-Text "55÷3=" "31÷9="
Replace-Text "12÷7=" "72÷5="
Replace-Text "52÷4=" "40÷3="
Replace-Text "43÷4=" "34÷9="

Replace-Text "24÷9=" "50÷5="
Replace-Text "29÷8=" "21÷4="
Replace-Text "20÷8=" "49÷4="
Replace-Text "73÷7=" "15÷5="
Replace-Text "57÷4=" "79÷8="

Replace-Text "81÷5=" "60÷4="
Replace-Text "96÷6=" "75÷6="
Replace-Text "78÷2=" "36÷4="
Replace-Text "89÷6=" "96÷2="
Replace-Text "15÷4=" "85÷8="

Replace-Text "19÷4=" "43÷2="
Replace-Text "64÷6=" "80÷7="
Replace-Text "88÷2=" "22÷8="
Replace-Text "66÷4=" "56÷2="
Replace-Text "26÷9=" "38÷6="
